# Update crypto price/volume data to reflect the latest scrape
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = "'58.784.12"
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = "'  +0.15%  "
$ws.Range('E2').Style = 'Normal'
$ws.Range('D3').Value = "'2.572.07"
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = "'  -0.64%  "
$ws.Range('E3').Style = 'Normal'
$ws.Range('E4').Value = "'  -0.04%  "
$ws.Range('E4').Style = 'Normal'
$ws.Range('D5').Value = "'562.95"
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = "'  +1.99%  "
$ws.Range('E5').Style = 'Normal'
$ws.Range('D6').Value = "'142.95"
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = "'  -0.94%  "
$ws.Range('E6').Style = 'Normal'
$ws.Range('D7').Value = "'0.999"
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = "'  +0.06%  "
$ws.Range('E7').Style = 'Normal'
$ws.Range('E8').Value = "'  +0.06%  "
$ws.Range('E8').Style = 'Normal'
$ws.Range('D9').Value = "'2.572.92"
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = "'  -0.92%  "
$ws.Range('E9').Style = 'Normal'
$ws.Range('E10').Value = "'  -2.34%  "
$ws.Range('E10').Style = 'Normal'
$ws.Range('E11').Value = "'  +2.04%  "
$ws.Range('E11').Style = 'Normal'
$ws.Range('D12').Value = "'0.152"
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = "'  +8.51%  "
$ws.Range('E12').Style = 'Normal'
$ws.Range('E13').Value = "'  +2.27%  "
$ws.Range('E13').Style = 'Normal'
$ws.Range('D14').Value = "'3.024.40"
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = "'  -0.57%  "
$ws.Range('E14').Style = 'Normal'
$ws.Range('D15').Value = "'58.861.57"
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = "'  +0.39%  "
$ws.Range('E15').Style = 'Normal'
$ws.Range('D16').Value = "'22.09"
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = "'  +6.52%  "
$ws.Range('E16').Style = 'Normal'
$ws.Range('E17').Value = "'  +3.66%  "
$ws.Range('E17').Style = 'Normal'
$ws.Range('D18').Value = "'2.571.00"
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = "'  -1.14%  "
$ws.Range('E18').Style = 'Normal'
$ws.Range('D20').Value = "'334.81"
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = "'  -0.48%  "
$ws.Range('E20').Style = 'Normal'
$ws.Range('E22').Value = "'  +0.41%  "
$ws.Range('E22').Style = 'Normal'
$ws.Range('E23').Value = "'  -0.01%  "
$ws.Range('E23').Style = 'Normal'
$ws.Range('D24').Value = "'63.92"
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = "'  -3.89%  "
$ws.Range('E24').Style = 'Normal'
$ws.Range('E25').Value = "'  +6.32%  "
$ws.Range('E25').Style = 'Normal'
$ws.Range('D26').Value = "'0.999"
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = "'  +0.12%  "
$ws.Range('E26').Style = 'Normal'
$ws.Range('E27').Value = "'  +1.83%  "
$ws.Range('E27').Style = 'Normal'
$ws.Range('D28').Value = "'7.22"
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = "'  +2.13%  "
$ws.Range('E28').Style = 'Normal'
$ws.Range('E29').Value = "'  +3.41%  "
$ws.Range('E29').Style = 'Normal'
$ws.Range('E30').Value = "'  +0.04%  "
$ws.Range('E30').Style = 'Normal'
$ws.Range('E31').Value = "'  -0.21%  "
$ws.Range('E31').Style = 'Normal'
$ws.Range('B32').Value = "'Aptos"
$ws.Range('B32').Style = 'Normal'
$ws.Range('C32').Value = "'https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range('C32').Style = 'Normal'
$ws.Range('D32').Value = "'6.02"
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = "'  +1.84%  "
$ws.Range('E32').Style = 'Normal'
$ws.Range('B33').Value = "'Monero"
$ws.Range('B33').Style = 'Normal'
$ws.Range('C33').Value = "'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range('C33').Style = 'Normal'
$ws.Range('D33').Value = "'157.93"
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = "'  +2.61%  "
$ws.Range('E33').Style = 'Normal'
$ws.Range('E34').Value = "'  +0.60%  "
$ws.Range('E34').Style = 'Normal'
$ws.Range('E35').Value = "'  +2.15%  "
$ws.Range('E35').Style = 'Normal'
$ws.Range('E36').Value = "'  +1.07%  "
$ws.Range('E36').Style = 'Normal'
$ws.Range('E37').Value = "'  +6.47%  "
$ws.Range('E37').Style = 'Normal'
$ws.Range('E38').Value = "'  +1.91%  "
$ws.Range('E38').Style = 'Normal'
$ws.Range('D39').Value = "'36.71"
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = "'  -1.37%  "
$ws.Range('E39').Style = 'Normal'
$ws.Range('D40').Value = "'1.49"
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = "'  +2.93%  "
$ws.Range('E40').Style = 'Normal'
$ws.Range('B41').Value = "'Filecoin"
$ws.Range('B41').Style = 'Normal'
$ws.Range('C41').Value = "'https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range('C41').Style = 'Normal'
$ws.Range('D41').Value = "'3.63"
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = "'  +1.07%  "
$ws.Range('E41').Style = 'Normal'
$ws.Range('B42').Value = "'Bittensor"
$ws.Range('B42').Style = 'Normal'
$ws.Range('C42').Value = "'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range('C42').Style = 'Normal'
$ws.Range('D42').Value = "'289.44"
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = "'  +3.18%  "
$ws.Range('E42').Style = 'Normal'
$ws.Range('E43').Value = "'  +0.16%  "
$ws.Range('E43').Style = 'Normal'
$ws.Range('E44').Value = "'  +2.16%  "
$ws.Range('E44').Style = 'Normal'
$ws.Range('D45').Value = "'0.594"
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = "'  -0.15%  "
$ws.Range('E45').Style = 'Normal'
$ws.Range('E46').Value = "'  -0.30%  "
$ws.Range('E46').Style = 'Normal'
$ws.Range('E47').Value = "'  +0.58%  "
$ws.Range('E47').Style = 'Normal'
$ws.Range('D48').Value = "'19.04"
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = "'  +1.86%  "
$ws.Range('E48').Style = 'Normal'
$ws.Range('D49').Value = "'123.90"
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = "'  +8.84%  "
$ws.Range('E49').Style = 'Normal'
$ws.Range('D50').Value = "'0.0231"
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = "'  +2.10%  "
$ws.Range('E50').Style = 'Normal'
$ws.Range('D51').Value = "'18.46"
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = "'  +3.59%  "
$ws.Range('E51').Style = 'Normal'
